$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mapping of row -> new DAMSLTag (col I) / DialogAct (col J) values
# Applies the dialog-act re-annotation from the SGNN re-run as per commit message.
$updates = @(
    @{Row=30; I='aa'; J='Agree/Accept'}
    @{Row=35; I='sv'; J='Statement-opinion'}
    @{Row=44; I='sd'; J='Statement-non-opinion'}
    @{Row=46; I='aa'; J='Agree/Accept'}
    @{Row=52; I='sv'; J='Statement-opinion'}
    @{Row=77; I='ba'; J='Appreciation'}
    @{Row=78; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=82; I='sv'; J='Statement-opinion'}
    @{Row=89; I='sv'; J='Statement-opinion'}
    @{Row=110; I='sd'; J='Statement-non-opinion'}
    @{Row=113; I='sv'; J='Statement-opinion'}
    @{Row=115; I='aa'; J='Agree/Accept'}
    @{Row=127; I='sv'; J='Statement-opinion'}
    @{Row=129; I='aa'; J='Agree/Accept'}
    @{Row=131; I='sd'; J='Statement-non-opinion'}
    @{Row=135; I='sd'; J='Statement-non-opinion'}
    @{Row=136; I='ba'; J='Appreciation'}
    @{Row=139; I='sd'; J='Statement-non-opinion'}
    @{Row=142; I='aa'; J='Agree/Accept'}
    @{Row=151; I='%'; J='Uninterpretable'}
    @{Row=169; I='aa'; J='Agree/Accept'}
    @{Row=173; I='aa'; J='Agree/Accept'}
    @{Row=176; I='sd'; J='Statement-non-opinion'}
    @{Row=178; I='sd'; J='Statement-non-opinion'}
    @{Row=206; I='sd'; J='Statement-non-opinion'}
    @{Row=209; I='aa'; J='Agree/Accept'}
    @{Row=212; I='sd'; J='Statement-non-opinion'}
    @{Row=215; I='aa'; J='Agree/Accept'}
    @{Row=232; I='sv'; J='Statement-opinion'}
    @{Row=239; I='sv'; J='Statement-opinion'}
    @{Row=243; I='sv'; J='Statement-opinion'}
    @{Row=244; I='sv'; J='Statement-opinion'}
    @{Row=247; I='sd'; J='Statement-non-opinion'}
    @{Row=250; I='sv'; J='Statement-opinion'}
    @{Row=262; I='sd'; J='Statement-non-opinion'}
    @{Row=295; I='sd'; J='Statement-non-opinion'}
    @{Row=311; I='sd'; J='Statement-non-opinion'}
    @{Row=323; I='sv'; J='Statement-opinion'}
    @{Row=325; I='sv'; J='Statement-opinion'}
    @{Row=340; I='sd'; J='Statement-non-opinion'}
    @{Row=341; I='aa'; J='Agree/Accept'}
    @{Row=342; I='aa'; J='Agree/Accept'}
    @{Row=358; I='sd'; J='Statement-non-opinion'}
    @{Row=361; I='sd'; J='Statement-non-opinion'}
    @{Row=373; I='sd'; J='Statement-non-opinion'}
    @{Row=380; I='sd'; J='Statement-non-opinion'}
    @{Row=385; I='sd'; J='Statement-non-opinion'}
    @{Row=398; I='sv'; J='Statement-opinion'}
    @{Row=406; I='sd'; J='Statement-non-opinion'}
    @{Row=410; I='sv'; J='Statement-opinion'}
    @{Row=413; I='sd'; J='Statement-non-opinion'}
    @{Row=421; I='sd'; J='Statement-non-opinion'}
    @{Row=426; I='sv'; J='Statement-opinion'}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}
